$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New monthly_mean (D) and weekly_share (F) values, recomputed 2020-10-23:
# monthly_mean is now the per-week average (old total / # weeks in month)
# weekly_share is recomputed as weekly_mean / (new) monthly_mean
$updates = @{
    2 = @(29.092, 0.917)
    3 = @(29.092, 1.06)
    4 = @(29.092, 0.989)
    5 = @(29.092, 1.05)
    6 = @(29.092, 0.985)
    7 = @(25.09, 1.057)
    8 = @(25.09, 1.016)
    9 = @(25.09, 0.969)
    10 = @(25.09, 0.921)
    11 = @(25.09, 1.037)
    12 = @(26.379, 1.021)
    13 = @(26.379, 1.013)
    14 = @(26.379, 0.97)
    15 = @(26.379, 0.995)
    16 = @(30.394, 0.903)
    17 = @(30.394, 0.981)
    18 = @(30.394, 1.076)
    19 = @(30.394, 1.039)
    20 = @(28.518, 1.153)
    21 = @(28.518, 1.058)
    22 = @(28.518, 0.97)
    23 = @(28.518, 0.917)
    24 = @(28.518, 0.901)
    25 = @(24.875, 1.123)
    26 = @(24.875, 0.864)
    27 = @(24.875, 0.888)
    28 = @(24.875, 1.125)
    29 = @(24.694, 1.091)
    30 = @(24.694, 0.943)
    31 = @(24.694, 0.974)
    32 = @(24.694, 0.992)
    33 = @(20.351, 1.101)
    34 = @(20.351, 1.04)
    35 = @(20.351, 0.923)
    36 = @(20.351, 1.071)
    37 = @(20.351, 0.865)
    38 = @(14.975, 1.057)
    39 = @(14.975, 0.973)
    40 = @(14.975, 0.955)
    41 = @(14.975, 1.015)
    42 = @(20.89, 0.821)
    43 = @(20.89, 0.82)
    44 = @(20.89, 0.903)
    45 = @(20.89, 1.264)
    46 = @(20.89, 1.192)
    47 = @(25.39, 1.064)
    48 = @(25.39, 0.954)
    49 = @(25.39, 1.024)
    50 = @(25.39, 0.958)
    51 = @(25.057, 0.967)
    52 = @(25.057, 0.937)
    53 = @(25.057, 1.06)
    54 = @(25.057, 1.036)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # column D: monthly_mean
    $ws.Cells.Item($row, 6).Value = $vals[1]   # column F: weekly_share
}

Write-Output "Updated monthly_mean and weekly_share for rows 2-54"
